$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "87.758.78"
$ws.Range("E2").Value = "  +3.43%  "
$ws.Range("D3").Value = "3.264.71"
$ws.Range("E3").Value = "  -1.49%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "'212.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.28%  "
$ws.Range("D6").Value = "'627.52"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.24%  "
$ws.Range("D7").Value = "'0.380"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +16.26%  "
$ws.Range("D8").Value = "'0.689"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +16.21%  "
$ws.Range("D9").Value = "'1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.14%  "
$ws.Range("D10").Value = "3.262.02"
$ws.Range("E10").Value = "  -1.26%  "
$ws.Range("E11").Value = "  -2.96%  "
$ws.Range("D12").Value = "'0.182"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +9.56%  "
$ws.Range("D13").Value = "'0.0000260"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -7.35%  "
$ws.Range("D14").Value = "'34.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.72%  "
$ws.Range("D15").Value = "3.870.21"
$ws.Range("E15").Value = "  -0.94%  "
$ws.Range("E16").Value = "  -1.88%  "
$ws.Range("D17").Value = "87.293.25"
$ws.Range("E17").Value = "  +3.13%  "
$ws.Range("D18").Value = "3.268.72"
$ws.Range("E18").Value = "  -0.31%  "
$ws.Range("E19").Value = "  -2.98%  "
$ws.Range("E20").Value = "  -3.97%  "
$ws.Range("D21").Value = "'434.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.00%  "
$ws.Range("E22").Value = "  -3.36%  "
$ws.Range("D23").Value = "'5.34"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.95%  "
$ws.Range("D24").Value = "'7.32"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.06%  "
$ws.Range("D25").Value = "'12.43"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.61%  "
$ws.Range("E26").Value = "  -6.53%  "
$ws.Range("D27").Value = "3.383.08"
$ws.Range("E27").Value = "  -1.96%  "
$ws.Range("D28").Value = "'76.59"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.97%  "
$ws.Range("E29").Value = "  -0.70%  "
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("D31").Value = "'0.182"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +12.24%  "
$ws.Range("D32").Value = "'0.999"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.14%  "
$ws.Range("D33").Value = "'8.79"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.45%  "
$ws.Range("D34").Value = "'550.75"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -8.19%  "
$ws.Range("D35").Value = "'1.40"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -12.34%  "
$ws.Range("E36").Value = "  -4.19%  "
$ws.Range("D37").Value = "'6.99"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +8.43%  "
$ws.Range("D38").Value = "'0.138"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -10.28%  "
$ws.Range("E39").Value = "  -3.51%  "
$ws.Range("D40").Value = "'1.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.22%  "
$ws.Range("D41").Value = "'21.73"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.63%  "
$ws.Range("B42").Value = "PolygonEcosystemToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D42").Value = "'0.393"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.64%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "'2.01"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.26%  "
$ws.Range("E44").Value = "  -6.00%  "
$ws.Range("E45").Value = "  +0.06%  "
$ws.Range("D46").Value = "'154.88"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.73%  "
$ws.Range("D47").Value = "'179.70"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.91%  "
$ws.Range("D48").Value = "'44.79"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.90%  "
$ws.Range("E49").Value = "  -4.22%  "
$ws.Range("D50").Value = "'4.23"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.39%  "
$ws.Range("D51").Value = "'0.123"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +10.63%  "
